# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Membrillo" (Vega Modelo de Temuco)
# above the current row 199, shifting the existing rows 199:213 down to 200:214.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 199:213 down by one row, creating a blank row 199.
$ws.Rows.Item(199).Insert()

# Fill the new row 199 with the latest weekly observation.
$ws.Range("A199").Value = 10
$ws.Range("B199").Value = "Vega Modelo de Temuco"
$ws.Range("C199").Value = "La Araucanía"
$ws.Range("D199").Value = 44783
$ws.Range("E199").Value = 9
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100104
$ws.Range("H199").Value = "Frutos de pepita"
$ws.Range("I199").Value = 100104003
$ws.Range("J199").Value = "Membrillo"
$ws.Range("K199").Value = "Champion"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 225
$ws.Range("N199").Value = 10000
$ws.Range("O199").Value = 10000
$ws.Range("P199").Value = 10000
$ws.Range("Q199").Value = "$/bandeja 18 kilos granel"
$ws.Range("R199").Value = "Región de O'Higgins"
$ws.Range("S199").Value = 556
$ws.Range("T199").Value = 18
